$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Rename the header row: "_old" -> "_FV2304", "_new" -> "_FV2310"
#    (column K stays "diff")
# ---------------------------------------------------------------------------
$headerRng = $ws.Range("A1:U1")

$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# ---------------------------------------------------------------------------
# 2. Turn the range A1:U67 into an Excel Table ("Table1"), preserving the
#    header row's existing cell formatting exactly (stash it on a scratch
#    row, strip the header style so the new table doesn't snapshot a
#    header dxf / table style name, then restore the formatting).
# ---------------------------------------------------------------------------
$scratchRng = $ws.Range("A69:U69")
$headerRng.Copy()
$scratchRng.PasteSpecial(-4122)   # xlPasteFormats

$headerRng.Style = "Normal"

$rng = $ws.Range("A1:U67")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"

$scratchRng.Copy()
$headerRng.PasteSpecial(-4122)    # xlPasteFormats
$scratchRng.Clear()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Freeze the header row (row 1).
# ---------------------------------------------------------------------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select()

Write-Host "Applied header rename, Table1 and frozen header row"
